$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.804.67"
$ws.Range("E2").Value = "  +2.25%  "
$ws.Range("D3").Value = "1.860.68"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6413"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.61%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.16"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3004"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07513"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "24.24"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07697"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "1.873.67"
$ws.Range("E13").Value = "  +2.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.069"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6870"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "84.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009406"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.096"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.58%  "
$ws.Range("D19").Value = "29.791.61"
$ws.Range("E19").Value = "  +2.24%  "
$ws.Range("D20").Value = "2.128.68"
$ws.Range("E20").Value = "  +2.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "241.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.55%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.463"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.52%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1432"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.580"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.19%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06109"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.47%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.506"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.273"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.171"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.133"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.873"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.159"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7340"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.610"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.863"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01802"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("D41").Value = "1.220.92"
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9322"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.284"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("D44").Value = "2.037.66"
$ws.Range("E44").Value = "  +3.37%  "
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "66.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000123"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5082"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.324"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4096"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.25%  "
